$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must stay as text
# (Excel would otherwise auto-convert them to numbers on assignment).
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated crypto price / volume data.
$ws.Range("D2").Value = "26.209.90"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.654.71"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "219.21"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "0.5244"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "0.2665"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "0.06361"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "20.68"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "0.07717"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.716.49"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.594"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "1.883.19"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "0.5633"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "0.0₅8257"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "65.39"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "26.215.36"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "4.699"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "10.39"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "191.97"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").Value = "6.005"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "143.75"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "0.1205"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("D27").Value = "7.272"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "1.514"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "0.05628"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").Value = "1.277"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "3.504"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").Value = "1.584"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "0.9534"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "2.800"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "2.413"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "0.5760"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "0.01599"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "6.010"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "0.8421"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").Value = "101.92"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "1.010.76"
$ws.Range("E44").Value = "  -6.20%  "
$ws.Range("D45").Value = "1.794.43"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "58.43"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "0.05343"
$ws.Range("E48").Value = "  +3.73%  "
$ws.Range("D49").Value = "8.027"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.09750"
$ws.Range("E51").Value = "  +1.71%  "

# Restore the default cell style now that the text is committed,
# so these cells do not keep a lingering custom number format.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
